# Change the table style ("Table_0" built-in style) applied to the
# finance-sources table on slide 6 to a different built-in table style.
#
#   old: {BB8EC9A1-7E65-4A46-B9CD-E6782DBD02BA}
#   new: {15E1BB81-DF1C-47CF-9471-16CB341D78BA}
#
# Table styles can't be assigned through the Style property directly;
# PowerPoint's object model requires Table.ApplyStyle("{GUID}").

$p = $ppt.ActivePresentation

$oldStyleId = "{BB8EC9A1-7E65-4A46-B9CD-E6782DBD02BA}"
$newStyleId = "{15E1BB81-DF1C-47CF-9471-16CB341D78BA}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
